# Swap the data (all columns except A, the row/id counter) between two
# row-pairs, matching the commit's "Atualização de bases das ligas" diff:
#   - row 195  <->  row 196
#   - row 210  <->  row 211

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB, $firstCol, $lastCol) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)

        $valA = $cellA.Value2
        $valB = $cellB.Value2

        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}

# Columns B (2) through AC (29) hold the match data; column A (1) is the
# sequential row id and must stay untouched.
Swap-Rows 195 196 2 29
Swap-Rows 210 211 2 29
